# Apply "Initial Glider Design v5.0" results update to Sheet1.
# Column A holds variable names, column B holds their numeric values.
# A handful of rows (41-46) were relabeled because two objective terms
# ("objective_sink", "objective_mass") were dropped and two new rows
# ("p_roll_max (rad/s)", "Cl_da (rad^-1)") were introduced in their place,
# shifting the remaining labels down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Relabel column A for rows whose variable name changed ---
$ws.Range("A41").Value = "p_roll_max (rad/s)"
$ws.Range("A42").Value = "t_roll (s)"
$ws.Range("A43").Value = "psi_0 (deg)"
$ws.Range("A44").Value = "Cl_da (rad^-1)"
$ws.Range("A45").Value = "objective_total"
$ws.Range("A46").Value = "objective_climb"

# --- Update recomputed numeric results in column B ---
# (Wrapped in [double]"..." so values in scientific notation, e.g. 1.2e-15,
#  parse correctly -- the PowerShell tokenizer here does not accept bare
#  "e-15"-style exponents as part of a numeric literal.)
$ws.Range("B2").Value = [double]"5.141061488811325"
$ws.Range("B3").Value = [double]"10.00000009984927"
$ws.Range("B4").Value = [double]"65.00000064741765"
$ws.Range("B5").Value = [double]"1.256345602990154"
$ws.Range("B6").Value = [double]"2.366201640490309"
$ws.Range("B7").Value = [double]"0.08277678496321403"
$ws.Range("B8").Value = [double]"8.087660995770321"
$ws.Range("B9").Value = [double]"1.504116959381735"
$ws.Range("B10").Value = [double]"1.41925860399273"
$ws.Range("B11").Value = [double]"-0.08485835538900477"
$ws.Range("B12").Value = [double]"0.03999999052013917"
$ws.Range("B13").Value = [double]"0.3999999900529829"
$ws.Range("B14").Value = [double]"0.01999999000524018"
$ws.Range("B15").Value = [double]"49463.23004659176"
$ws.Range("B16").Value = [double]"-0.1015359301217917"
$ws.Range("B17").Value = [double]"0.8056713410028052"
$ws.Range("B18").Value = [double]"0.8114614932330131"
$ws.Range("B19").Value = [double]"2.053750223110356"
$ws.Range("B20").Value = [double]"0.1405388941750224"
$ws.Range("B22").Value = [double]"0.1140419009245801"
$ws.Range("B23").Value = [double]"0.8056713410028052"
$ws.Range("B24").Value = [double]"0.1784068470692191"
$ws.Range("B25").Value = [double]"0.04460171176730478"
$ws.Range("B27").Value = [double]"0.007957250770294935"
$ws.Range("B28").Value = [double]"0.8056713410028052"
$ws.Range("B29").Value = [double]"0.0677824264146713"
$ws.Range("B30").Value = [double]"0.03389121320733565"
$ws.Range("B32").Value = [double]"0.002297228665330165"
$ws.Range("B36").Value = [double]"-1.32197450498956e-15"
$ws.Range("B37").Value = [double]"1.608227095629997e-14"
$ws.Range("B38").Value = [double]"10.00580694562281"
$ws.Range("B39").Value = [double]"1.055098752613763"
$ws.Range("B40").Value = [double]"60.45270549428633"
$ws.Range("B41").Value = [double]"13.61794718233036"
$ws.Range("B42").Value = [double]"1.075220705441564"
$ws.Range("B43").Value = [double]"-31.60949773307858"
$ws.Range("B44").Value = [double]"-0.4866281482559692"
$ws.Range("B45").Value = [double]"1.065901365071977"
$ws.Range("B46").Value = [double]"0.007249436628742574"
$ws.Range("B47").Value = [double]"1.057650766716904"
$ws.Range("B48").Value = [double]"0.001001161726330736"
$ws.Range("B49").Value = [double]"3.333585406817299e-10"
$ws.Range("B50").Value = [double]"0.03284266895196183"
$ws.Range("B51").Value = [double]"-1.27227663362964e-17"
$ws.Range("B52").Value = [double]"0.003182731310932615"
$ws.Range("B53").Value = [double]"0.001945924981276573"
$ws.Range("B54").Value = [double]"0.001961788971508408"
$ws.Range("B55").Value = [double]"0.003903935294542807"
$ws.Range("B56").Value = [double]"1.387152435409648e-21"
$ws.Range("B57").Value = [double]"-4.48775121766283e-06"
$ws.Range("B58").Value = [double]"3.399195133571679e-21"
$ws.Range("B59").Value = [double]"0.03289764424912821"
$ws.Range("B60").Value = [double]"0.0006418656542902716"
$ws.Range("B61").Value = [double]"0.0001408059179774253"
$ws.Range("B66").Value = [double]"0.008164865440121372"
$ws.Range("B68").Value = [double]"-9.99928206569716e-09"
$ws.Range("B69").Value = [double]"0.006131613700978818"
$ws.Range("B70").Value = [double]"-0.0351347235437556"
$ws.Range("B73").Value = [double]"-0.0351347235437556"
$ws.Range("B74").Value = [double]"0.4054701244165568"
$ws.Range("B75").Value = [double]"0.01454018416698145"
$ws.Range("B76").Value = [double]"0.7610696292355005"
$ws.Range("B79").Value = [double]"0.7610696292355005"
$ws.Range("B80").Value = [double]"0.08920342353460955"
$ws.Range("B82").Value = [double]"0.7717801277954696"
$ws.Range("B85").Value = [double]"0.7717801277954696"
$ws.Range("B87").Value = [double]"0.0677824264146713"
$ws.Range("B88").Value = [double]"-0.09968772864393802"
$ws.Range("B89").Value = [double]"-1.642741440261673e-18"
$ws.Range("B90").Value = [double]"1.933514829369654"
$ws.Range("B91").Value = [double]"-8.192865265399014e-17"
$ws.Range("B92").Value = [double]"2.283188828972049e-14"
$ws.Range("B93").Value = [double]"-2.781845244806068e-17"
$ws.Range("B94").Value = [double]"0.09968772864393802"
$ws.Range("B95").Value = [double]"-0.2375780819118086"
$ws.Range("B96").Value = [double]"8.192865265399014e-17"
$ws.Range("B97").Value = [double]"8.551459592356331e-17"
$ws.Range("B98").Value = [double]"1.92145098652806"
$ws.Range("B99").Value = [double]"-1.642741440261673e-18"
$ws.Range("B100").Value = [double]"0.2375780819118086"
$ws.Range("B101").Value = [double]"8.192865265399014e-17"
$ws.Range("B102").Value = [double]"2.283188828972049e-14"
$ws.Range("B103").Value = [double]"2.781845244806068e-17"
$ws.Range("B104").Value = [double]"1.040766902446284"
$ws.Range("B105").Value = [double]"-8.898019945804752e-19"
$ws.Range("B106").Value = [double]"0.1286857724366271"
$ws.Range("B107").Value = [double]"5.468800225985341e-17"
$ws.Range("B108").Value = [double]"8.799732308183184e-14"
$ws.Range("B109").Value = [double]"1.85690297724074e-17"
$ws.Range("B110").Value = [double]"1.854331488203392"
$ws.Range("B111").Value = [double]"0.06870129053957613"
$ws.Range("B112").Value = [double]"1.214306433183765e-17"
$ws.Range("B113").Value = [double]"8.326672684688674e-17"
$ws.Range("B114").Value = [double]"0.05210876439358686"
$ws.Range("B115").Value = [double]"1.734723475976807e-17"
$ws.Range("B116").Value = [double]"0.8109404162444613"
$ws.Range("B117").Value = [double]"0.7998638959219753"
$ws.Range("B118").Value = [double]"0.06518784432553081"
$ws.Range("B119").Value = [double]"0.02827084386385509"
$ws.Range("B121").Value = [double]"-8.944667923005412e-19"
$ws.Range("B122").Value = [double]"-0.0513365394690181"
$ws.Range("B123").Value = [double]"2.236166980751353e-19"
$ws.Range("B124").Value = [double]"0.1784068470692191"
$ws.Range("B126").Value = [double]"2.710505431213761e-20"
$ws.Range("B127").Value = [double]"0.001118499803066326"
$ws.Range("B128").Value = [double]"-1.378580577209932e-17"
$ws.Range("B129").Value = [double]"-4.161356364208345e-19"
$ws.Range("B130").Value = [double]"-0.0001109254283730531"
$ws.Range("B131").Value = [double]"1.023332831614988e-17"
$ws.Range("B132").Value = [double]"0.0677824264146713"
$ws.Range("B134").Value = [double]"0.001931653999137287"
$ws.Range("B135").Value = [double]"0.001479730226077065"
$ws.Range("B137").Value = [double]"-2.747176417521673e-20"
$ws.Range("B138").Value = [double]"-0.0006612994961728788"
$ws.Range("B139").Value = [double]"1.427267406759412e-20"
$ws.Range("B142").Value = [double]"0.09957036443257461"
$ws.Range("B143").Value = [double]"0.138007717479234"
$ws.Range("B144").Value = [double]"4.021942757649293"
$ws.Range("B145").Value = [double]"0.7196226212826305"
$ws.Range("B146").Value = [double]"2.785088802896838e-22"
$ws.Range("B147").Value = [double]"5.286877008665689e-13"
$ws.Range("B148").Value = [double]"-0.1608776721785141"
$ws.Range("B149").Value = [double]"-2.325967398586828e-13"
$ws.Range("B150").Value = [double]"0.03846422338667357"
$ws.Range("B151").Value = [double]"-5.438646727552807e-06"
$ws.Range("B152").Value = [double]"-6.658968308846092e-07"
$ws.Range("B153").Value = [double]"0.1396629135701793"
$ws.Range("B154").Value = [double]"-0.02499999015796307"
$ws.Range("B155").Value = [double]"6.799741356075048e-07"
$ws.Range("B156").Value = [double]"0.03217912956944434"
$ws.Range("B157").Value = [double]"-0.1541226739859743"
$ws.Range("B158").Value = [double]"-0.00824421885603499"
$ws.Range("B159").Value = [double]"-0.001687347274764717"
$ws.Range("B160").Value = [double]"-0.07811761689330377"
$ws.Range("B161").Value = [double]"-0.6462573293282363"
$ws.Range("B162").Value = [double]"-0.001214576490396257"
$ws.Range("B163").Value = [double]"-0.2729271052282198"
$ws.Range("B164").Value = [double]"0.7477096696717478"
$ws.Range("B165").Value = [double]"1.261830316577978"
$ws.Range("B166").Value = [double]"-5.659664324254516e-15"
$ws.Range("B167").Value = [double]"-5.102221748878642e-14"
$ws.Range("B168").Value = [double]"-16.67366635714174"
$ws.Range("B169").Value = [double]"-1.155931116872107e-14"
$ws.Range("B170").Value = [double]"-0.0001276918786263792"
$ws.Range("B171").Value = [double]"-3.21477403031345e-05"
$ws.Range("B172").Value = [double]"0.1128903331499223"
$ws.Range("B173").Value = [double]"0.2225041885750886"
$ws.Range("B174").Value = [double]"-0.0002291904827817157"
$ws.Range("B175").Value = [double]"-0.07231340751734523"
